# Mark specific vocabulary rows as "Processed" in column C.
# Rows 21-25 plus every row ending in 5 from row 45 through 1035
# (45, 55, 65, ... 1035) get a new C cell with the text "Processed".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(21,22,23,24,25) + @(45..1035 | Where-Object { ($_ - 45) % 10 -eq 0 })

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
